$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H69").Value = 4727.087
$ws.Range("I69").Value = 4367
$ws.Range("J69").Value = 4958.5713
$ws.Range("K69").Value = 13101
$ws.Range("L69").Value = 14875.7139
$ws.Range("M69").Value = -12227
$ws.Range("N69").Value = -16623.7139

$ws.Range("H72").Value = 4727.087
$ws.Range("I72").Value = 4367
$ws.Range("J72").Value = 4958.5713
$ws.Range("K72").Value = 39303
$ws.Range("L72").Value = 44627.14169999999
$ws.Range("M72").Value = -34935
$ws.Range("N72").Value = -53363.14169999999

$ws.Range("H76").Value = 7138.143
$ws.Range("I76").Value = 10150.429
$ws.Range("K76").Value = 10150.429
$ws.Range("M76").Value = -9835.429

$ws.Range("H79").Value = 7138.143
$ws.Range("I79").Value = 10150.429
$ws.Range("K79").Value = 10150.429
$ws.Range("M79").Value = -9058.429

$ws.Range("H98").Value = 7333.6665
$ws.Range("I98").Value = 7333.6665
$ws.Range("K98").Value = 7333.6665
$ws.Range("M98").Value = -5835.6665

$ws.Range("H122").Value = 7333.6665
$ws.Range("I122").Value = 7333.6665
$ws.Range("K122").Value = 22000.9995
$ws.Range("M122").Value = -19550.9995

$ws.Range("H129").Value = 894.75757
$ws.Range("I129").Value = 472.3
$ws.Range("J129").Value = 942.22473
$ws.Range("K129").Value = 1416.9
$ws.Range("L129").Value = 2826.67419
$ws.Range("M129").Value = 3583.1
$ws.Range("N129").Value = -12826.67419

$ws.Range("H137").Value = 1245.0605
$ws.Range("I137").Value = 1039.931
$ws.Range("J137").Value = 2732.25
$ws.Range("K137").Value = 3119.793
$ws.Range("L137").Value = 8196.75
$ws.Range("M137").Value = -569.7930000000001
$ws.Range("N137").Value = -13296.75

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 5225.6895
$ws.Range("I32").Value = 4623.391
$ws.Range("J32").Value = 7534.5
$ws.Range("K32").Value = 4623.391
$ws.Range("L32").Value = 7534.5
$ws.Range("M32").Value = -4336.391
$ws.Range("N32").Value = -8108.5

$ws.Range("H61").Value = 5959.64
$ws.Range("I61").Value = 6124.625
$ws.Range("J61").Value = 2000
$ws.Range("K61").Value = 6124.625
$ws.Range("L61").Value = 2000
$ws.Range("M61").Value = -5912.625
$ws.Range("N61").Value = -2424

$ws.Range("H74").Value = 2127.5833
$ws.Range("I74").Value = 1886
$ws.Range("J74").Value = 2714.2856
$ws.Range("K74").Value = 1886
$ws.Range("L74").Value = 2714.2856
$ws.Range("M74").Value = -1012
$ws.Range("N74").Value = -4462.2856

$ws.Range("H77").Value = 2127.5833
$ws.Range("I77").Value = 1886
$ws.Range("J77").Value = 2714.2856
$ws.Range("K77").Value = 9430
$ws.Range("L77").Value = 13571.428
$ws.Range("M77").Value = -5062
$ws.Range("N77").Value = -22307.428

$ws.Range("H132").Value = 2970.0857
$ws.Range("I132").Value = 1225.5454
$ws.Range("K132").Value = 3676.6362
$ws.Range("M132").Value = -1146.6362

$ws.Range("H136").Value = 5959.64
$ws.Range("I136").Value = 6124.625
$ws.Range("J136").Value = 2000
$ws.Range("K136").Value = 18373.875
$ws.Range("L136").Value = 6000
$ws.Range("M136").Value = -15823.875
$ws.Range("N136").Value = -11100

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H50").Value = 40780
$ws.Range("J50").Value = 40780
$ws.Range("L50").Value = 40780
$ws.Range("N50").Value = -41928

$ws.Range("H134").Value = 3848.2932
$ws.Range("I134").Value = 4931.242
$ws.Range("K134").Value = 14793.726
$ws.Range("M134").Value = -12258.726

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 5244.5366
$ws.Range("I31").Value = 1350.0769
$ws.Range("J31").Value = 11994.934
$ws.Range("K31").Value = 1350.0769
$ws.Range("L31").Value = 11994.934
$ws.Range("M31").Value = -1055.0769
$ws.Range("N31").Value = -12584.934

$ws.Range("H34").Value = 5244.5366
$ws.Range("I34").Value = 1350.0769
$ws.Range("J34").Value = 11994.934
$ws.Range("K34").Value = 1350.0769
$ws.Range("L34").Value = 11994.934
$ws.Range("M34").Value = -1148.0769
$ws.Range("N34").Value = -12398.934

$ws.Range("H86").Value = 2376.7646
$ws.Range("I86").Value = 2515.5833
$ws.Range("J86").Value = 2043.6
$ws.Range("K86").Value = 2515.5833
$ws.Range("L86").Value = 2043.6
$ws.Range("M86").Value = -1392.5833
$ws.Range("N86").Value = -4289.6

$ws.Range("H89").Value = 2376.7646
$ws.Range("I89").Value = 2515.5833
$ws.Range("J89").Value = 2043.6
$ws.Range("K89").Value = 12577.9165
$ws.Range("L89").Value = 10218
$ws.Range("M89").Value = -6961.916499999999
$ws.Range("N89").Value = -21450

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H37").Value = 42100
$ws.Range("J37").Value = 42100
$ws.Range("L37").Value = 126300
$ws.Range("N37").Value = -126524

$ws.Range("H107").Value = 125655.94
$ws.Range("I107").Value = 230
$ws.Range("J107").Value = 154600.39
$ws.Range("K107").Value = 690
$ws.Range("L107").Value = 463801.17
$ws.Range("M107").Value = 1230
$ws.Range("N107").Value = -467641.17

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 6357.968
$ws.Range("I70").Value = 6472.636
$ws.Range("J70").Value = 6077.6665
$ws.Range("K70").Value = 6472.636
$ws.Range("L70").Value = 6077.6665
$ws.Range("M70").Value = -6202.636
$ws.Range("N70").Value = -6617.6665

$ws.Range("H73").Value = 6357.968
$ws.Range("I73").Value = 6472.636
$ws.Range("J73").Value = 6077.6665
$ws.Range("K73").Value = 6472.636
$ws.Range("L73").Value = 6077.6665
$ws.Range("M73").Value = -5536.636
$ws.Range("N73").Value = -7949.6665

$ws.Range("H122").Value = 3481120
$ws.Range("I122").Value = 2494786.8
$ws.Range("J122").Value = 7144643
$ws.Range("K122").Value = 7484360.399999999
$ws.Range("L122").Value = 21433929
$ws.Range("M122").Value = -7481910.399999999
$ws.Range("N122").Value = -21438829

$ws.Range("H132").Value = 3585.3872
$ws.Range("I132").Value = 3962.2778
$ws.Range("K132").Value = 11886.8334
$ws.Range("M132").Value = -9356.8334

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H68").Value = 37039020
$ws.Range("I68").Value = 1545.3334
$ws.Range("K68").Value = 1545.3334
$ws.Range("M68").Value = -796.3334

$ws.Range("H71").Value = 37039020
$ws.Range("I71").Value = 1545.3334
$ws.Range("K71").Value = 7726.666999999999
$ws.Range("M71").Value = -3982.666999999999

$ws.Range("H122").Value = 14288115
$ws.Range("I122").Value = 35715784
$ws.Range("J122").Value = 3001.6667
$ws.Range("K122").Value = 107147352
$ws.Range("L122").Value = 9005.000100000001
$ws.Range("M122").Value = -107144902
$ws.Range("N122").Value = -13905.0001

$ws.Range("H132").Value = 19935188
$ws.Range("I132").Value = 28655032
$ws.Range("J132").Value = 4113.2856
$ws.Range("K132").Value = 85965096
$ws.Range("L132").Value = 12339.8568
$ws.Range("M132").Value = -85962566
$ws.Range("N132").Value = -17399.8568

$ws.Range("H136").Value = 5887.6875
$ws.Range("I136").Value = 6121.643
$ws.Range("J136").Value = 4250
$ws.Range("K136").Value = 18364.929
$ws.Range("L136").Value = 12750
$ws.Range("M136").Value = -15814.929
$ws.Range("N136").Value = -17850

$ws.Range("H141").Value = 50383.332
$ws.Range("J141").Value = 50383.332
$ws.Range("L141").Value = 50383.332
$ws.Range("N141").Value = -60743.332

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 13514.714
$ws.Range("I62").Value = 3750
$ws.Range("J62").Value = 26534.334
$ws.Range("K62").Value = 3750
$ws.Range("L62").Value = 26534.334
$ws.Range("M62").Value = -3126
$ws.Range("N62").Value = -27782.334

$ws.Range("H65").Value = 13514.714
$ws.Range("I65").Value = 3750
$ws.Range("J65").Value = 26534.334
$ws.Range("K65").Value = 18750
$ws.Range("L65").Value = 132671.67
$ws.Range("M65").Value = -15630
$ws.Range("N65").Value = -138911.67

$ws.Range("H81").Value = 1892.6666
$ws.Range("I81").Value = 1506
$ws.Range("J81").Value = 2376
$ws.Range("K81").Value = 3012
$ws.Range("L81").Value = 4752
$ws.Range("M81").Value = -1951
$ws.Range("N81").Value = -6874

$ws.Range("H84").Value = 1892.6666
$ws.Range("I84").Value = 1506
$ws.Range("J84").Value = 2376
$ws.Range("K84").Value = 15060
$ws.Range("L84").Value = 23760
$ws.Range("M84").Value = -9756
$ws.Range("N84").Value = -34368

$ws.Range("H122").Value = 3654.889
$ws.Range("I122").Value = 2999.25
$ws.Range("J122").Value = 4179.4
$ws.Range("K122").Value = 8997.75
$ws.Range("L122").Value = 12538.2
$ws.Range("M122").Value = -6547.75
$ws.Range("N122").Value = -17438.2

$ws.Range("H132").Value = 2193.4167
$ws.Range("I132").Value = 1775.8182
$ws.Range("K132").Value = 5327.4546
$ws.Range("M132").Value = -2797.4546

$ws.Range("H136").Value = 2870.6206
$ws.Range("I136").Value = 3459.0557
$ws.Range("J136").Value = 1907.7273
$ws.Range("K136").Value = 10377.1671
$ws.Range("L136").Value = 5723.1819
$ws.Range("M136").Value = -7827.167099999999
$ws.Range("N136").Value = -10823.1819

$ws.Range("H137").Value = 45425
$ws.Range("J137").Value = 45425
$ws.Range("L137").Value = 45425
$ws.Range("N137").Value = -55625

$ws.Range("H139").Value = 79850
$ws.Range("J139").Value = 79850
$ws.Range("L139").Value = 79850
$ws.Range("N139").Value = -90130
